# Apply the NIT-9002325001 "Estado de Cuenta" update:
#  - Remove worker LUIS FERNANDO LEOTTAU FONSECA (CC 8834275) entirely (rows 16-17)
#  - Keep worker ADRIANA DEL SOCORRO JIMENEZ ALMEIDA (CC 32939680) but flip her period
#    rows into ascending order (1611 -> 2104) and refresh "Salario Basico" to 781242
#  - Refresh the summary header values (Valor Mora total, worker count, period count)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the first worker's two rows (periods 1606 and 1605) ---------------
$ws.Rows("16:17").Delete() | Out-Null

# --- 2. Capture the remaining worker's period (E) and "Valor Mora" (F) values ----
$firstDataRow = 16
$lastDataRow = 69

$periods = @()
$valores = @()
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $periods += , $ws.Cells.Item($r, 5).Value2
    $valores += , $ws.Cells.Item($r, 6).Value2
}

# --- 3. Write the values back in reverse row order (ascending period order) -----
$count = $periods.Count
for ($k = 0; $k -lt $count; $k++) {
    $r = $firstDataRow + $k
    $srcIndex = $count - 1 - $k
    $ws.Cells.Item($r, 5).Value = $periods[$srcIndex]
    $ws.Cells.Item($r, 6).Value = $valores[$srcIndex]
    $ws.Cells.Item($r, 7).Value = 781242
}

# --- 4. Refresh the summary header values ----------------------------------------
$ws.Range("E11").Value = 1602518
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 54
